# Update market-price / profit figures across sheets (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 9731.666999999999
$ws.Range("I64").Value = 6200
$ws.Range("J64").Value = 12086.111
$ws.Range("K64").Value = 6200
$ws.Range("L64").Value = 12086.111
$ws.Range("M64").Value = -5952
$ws.Range("N64").Value = -12582.111

$ws.Range("H67").Value = 9731.666999999999
$ws.Range("I67").Value = 6200
$ws.Range("J67").Value = 12086.111
$ws.Range("K67").Value = 6200
$ws.Range("L67").Value = 12086.111
$ws.Range("M67").Value = -5342
$ws.Range("N67").Value = -13802.111

$ws.Range("H98").Value = 3001.8462
$ws.Range("I98").Value = 3001.8462
$ws.Range("K98").Value = 3001.8462
$ws.Range("M98").Value = -1503.8462

$ws.Range("H100").Value = 1231.25
$ws.Range("I100").Value = 1242.5
$ws.Range("J100").Value = 1197.5
$ws.Range("K100").Value = 1242.5
$ws.Range("L100").Value = 1197.5
$ws.Range("M100").Value = -701.5
$ws.Range("N100").Value = -2279.5

$ws.Range("H122").Value = 3001.8462
$ws.Range("I122").Value = 3001.8462
$ws.Range("K122").Value = 9005.5386
$ws.Range("M122").Value = -6555.5386

$ws.Range("H137").Value = 49807.785
$ws.Range("I137").Value = 52869.92
$ws.Range("J137").Value = 10000
$ws.Range("K137").Value = 158609.76
$ws.Range("L137").Value = 30000
$ws.Range("M137").Value = -156059.76
$ws.Range("N137").Value = -35100

$ws.Range("H138").Value = 9829.666999999999
$ws.Range("J138").Value = 9836
$ws.Range("L138").Value = 29508
$ws.Range("N138").Value = -39788

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2844.4934
$ws.Range("I32").Value = 2132.554
$ws.Range("K32").Value = 2132.554
$ws.Range("M32").Value = -1845.554

$ws.Range("H61").Value = 7580730.5
$ws.Range("I61").Value = 10420736
$ws.Range("K61").Value = 10420736
$ws.Range("M61").Value = -10420524

$ws.Range("H74").Value = 32753.5
$ws.Range("I74").Value = 2197.8333
$ws.Range("K74").Value = 2197.8333
$ws.Range("M74").Value = -1323.8333

$ws.Range("H77").Value = 32753.5
$ws.Range("I77").Value = 2197.8333
$ws.Range("K77").Value = 10989.1665
$ws.Range("M77").Value = -6621.166499999999

$ws.Range("H88").Value = 1449
$ws.Range("I88").Value = 1598.3334
$ws.Range("J88").Value = 1299.6666
$ws.Range("K88").Value = 1598.3334
$ws.Range("L88").Value = 1299.6666
$ws.Range("M88").Value = -1192.3334
$ws.Range("N88").Value = -2111.6666

$ws.Range("H91").Value = 1449
$ws.Range("I91").Value = 1598.3334
$ws.Range("J91").Value = 1299.6666
$ws.Range("K91").Value = 1598.3334
$ws.Range("L91").Value = 1299.6666
$ws.Range("M91").Value = -194.3334
$ws.Range("N91").Value = -4107.6666

$ws.Range("H110").Value = 8216.333000000001
$ws.Range("I110").Value = 4799.4
$ws.Range("J110").Value = 12487.5
$ws.Range("K110").Value = 4799.4
$ws.Range("L110").Value = 12487.5
$ws.Range("M110").Value = -2754.4
$ws.Range("N110").Value = -16577.5

$ws.Range("H133").Value = 195199.6
$ws.Range("J133").Value = 194999.5
$ws.Range("L133").Value = 194999.5
$ws.Range("N133").Value = -200059.5

$ws.Range("H136").Value = 7580730.5
$ws.Range("I136").Value = 10420736
$ws.Range("K136").Value = 31262208
$ws.Range("M136").Value = -31259658

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4960
$ws.Range("I107").Value = 5496.857
$ws.Range("K107").Value = 5496.857
$ws.Range("M107").Value = -3576.857

$ws.Range("H108").Value = 38999
$ws.Range("J108").Value = 38999
$ws.Range("L108").Value = 38999
$ws.Range("N108").Value = -46679

$ws.Range("H134").Value = 3983.1035
$ws.Range("I134").Value = 3848.5715
$ws.Range("K134").Value = 11545.7145
$ws.Range("M134").Value = -9010.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 207.18182
$ws.Range("I7").Value = 97.5
$ws.Range("J7").Value = 499.66666
$ws.Range("K7").Value = 97.5
$ws.Range("L7").Value = 499.66666
$ws.Range("M7").Value = 15.5
$ws.Range("N7").Value = -725.66666

$ws.Range("H16").Value = 1801.091
$ws.Range("I16").Value = 1922.4
$ws.Range("J16").Value = 1541.1428
$ws.Range("K16").Value = 1922.4
$ws.Range("L16").Value = 1541.1428
$ws.Range("M16").Value = -1635.4
$ws.Range("N16").Value = -2115.1428

$ws.Range("H22").Value = 199.3125
$ws.Range("J22").Value = 239.66667
$ws.Range("L22").Value = 239.66667
$ws.Range("N22").Value = -939.6666700000001

$ws.Range("H35").Value = 4998
$ws.Range("I35").Value = 4998
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 4998
$ws.Range("M35").Value = -4704
$ws.Range("N35").ClearContents()

$ws.Range("H105").Value = 1029.2858
$ws.Range("I105").Value = 961.2
$ws.Range("K105").Value = 961.2
$ws.Range("M105").Value = 785.8

$ws.Range("H107").Value = 3543.6428
$ws.Range("I107").Value = 6352.75
$ws.Range("K107").Value = 6352.75
$ws.Range("M107").Value = -4432.75

$ws.Range("H113").Value = 1801.091
$ws.Range("I113").Value = 1922.4
$ws.Range("J113").Value = 1541.1428
$ws.Range("K113").Value = 1922.4
$ws.Range("L113").Value = 1541.1428
$ws.Range("M113").Value = 247.5999999999999
$ws.Range("N113").Value = -5881.1428

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 22224498
$ws.Range("I131").Value = 62500656
$ws.Range("J131").Value = 5266114.5
$ws.Range("K131").Value = 187501968
$ws.Range("L131").Value = 15798343.5
$ws.Range("M131").Value = -187496928
$ws.Range("N131").Value = -15808423.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 29997
$ws.Range("J33").Value = 29997
$ws.Range("L33").Value = 29997
$ws.Range("N33").Value = -30501

$ws.Range("H70").Value = 15027.111
$ws.Range("I70").Value = 11312.25
$ws.Range("J70").Value = 17999
$ws.Range("K70").Value = 11312.25
$ws.Range("L70").Value = 17999
$ws.Range("M70").Value = -11042.25
$ws.Range("N70").Value = -18539

$ws.Range("H73").Value = 15027.111
$ws.Range("I73").Value = 11312.25
$ws.Range("J73").Value = 17999
$ws.Range("K73").Value = 11312.25
$ws.Range("L73").Value = 17999
$ws.Range("M73").Value = -10376.25
$ws.Range("N73").Value = -19871

$ws.Range("H113").Value = 185321.45
$ws.Range("I113").Value = 214320.38
$ws.Range("J113").Value = 1661.6666
$ws.Range("K113").Value = 214320.38
$ws.Range("L113").Value = 1661.6666
$ws.Range("M113").Value = -212150.38
$ws.Range("N113").Value = -6001.6666

$ws.Range("H122").Value = 4286.125
$ws.Range("I122").Value = 5188.9
$ws.Range("K122").Value = 15566.7
$ws.Range("M122").Value = -13116.7

$ws.Range("H126").Value = 10702.2
$ws.Range("I126").Value = 10670.333
$ws.Range("K126").Value = 32010.999
$ws.Range("M126").Value = -29540.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2082.6667
$ws.Range("I22").Value = 1399.4
$ws.Range("J22").Value = 2936.75
$ws.Range("K22").Value = 1399.4
$ws.Range("L22").Value = 2936.75
$ws.Range("M22").Value = -1104.4
$ws.Range("N22").Value = -3526.75

$ws.Range("H27").Value = 2082.6667
$ws.Range("I27").Value = 1399.4
$ws.Range("J27").Value = 2936.75
$ws.Range("K27").Value = 1399.4
$ws.Range("L27").Value = 2936.75
$ws.Range("M27").Value = -1292.4
$ws.Range("N27").Value = -3150.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 23983.334
$ws.Range("J69").Value = 23983.334
$ws.Range("L69").Value = 23983.334
$ws.Range("N69").Value = -25481.334

$ws.Range("H72").Value = 23983.334
$ws.Range("J72").Value = 23983.334
$ws.Range("L72").Value = 71950.00199999999
$ws.Range("N72").Value = -79438.00199999999

Write-Host "Updated ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR profit figures"
